# Update countries & provincias Spain
# - Refresh "last updated" timestamp
# - Reorder a couple of shared-string rows so the alphabetised country list
#   stays correct (Belice before Birmania/Islas Feroe; Islas Malvinas before
#   Montserrat) and refresh the COVID case counters that changed between
#   pulls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a ..." timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 06:47"

# --- Kazajistan (row 29) ---------------------------------------------------
$ws.Range("D29").Value = 80716
$ws.Range("E29").Value = 20302

# --- Haiti (row 96) ---------------------------------------------------------
$ws.Range("B96").Value = 7831
$ws.Range("C96").Value = 21
$ws.Range("D96").Value = 5235
$ws.Range("E96").Value = 2400
$ws.Range("G96").Value = 4
$ws.Range("H96").Value = 196

# --- Tailandia (row 117) ----------------------------------------------------
$ws.Range("D117").Value = 3193
$ws.Range("E117").Value = 125

# --- Belice / Birmania / Islas Feroe reorder (rows 172-174) ---------------
# "Belice" moves up to sit right after "Comoras" (ahead of "Birmania" and
# "Islas Feroe"), each picking up the data the row above used to show, and
# "Belice" itself getting freshly updated counters.
$ws.Range("A172").Value = "Belice"
$ws.Range("B172").Value = 388
$ws.Range("C172").Value = 32
$ws.Range("D172").Value = 35
$ws.Range("E172").Value = 350
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = 3

$ws.Range("A173").Value = "Birmania"
$ws.Range("B173").Value = 374
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 322
$ws.Range("E173").Value = 46
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 6

$ws.Range("A174").Value = "Islas Feroe"
$ws.Range("B174").Value = 365
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 225
$ws.Range("E174").Value = 140
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# --- Mongolia (row 178) -----------------------------------------------------
$ws.Range("B178").Value = 298
$ws.Range("C178").Value = 1
$ws.Range("D178").Value = 272
$ws.Range("E178").Value = 26

# --- Camboya (row 181) ------------------------------------------------------
$ws.Range("D181").Value = 229
$ws.Range("E181").Value = 44

# --- Butan (row 191) --------------------------------------------------------
$ws.Range("B191").Value = 133
$ws.Range("C191").Value = 2
$ws.Range("D191").Value = 102

# --- Islas Malvinas / Montserrat reorder (rows 213-214) --------------------
# "Islas Malvinas" moves up ahead of "Montserrat"; data travels with the
# country name, values themselves are unchanged, only the order swaps.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
